$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the shared-strings table in the same order as the target workbook:
# Storage tank, Deliver Station, Pumping station.
$ws.Cells.Item(11, 2).Value = "Storage tank"
$ws.Cells.Item(17, 2).Value = "Deliver Station"
$ws.Cells.Item(2, 2).Value = "Pumping station"

# Rows 2-10: "Pump Stations" -> "Pumping station"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "Pumping station"
}

# Rows 11-16: "Storage Tanks" -> "Storage tank"
for ($r = 11; $r -le 16; $r++) {
    $ws.Cells.Item($r, 2).Value = "Storage tank"
}

# Rows 17-50: "Delivery Nodes" -> "Deliver Station"
for ($r = 17; $r -le 50; $r++) {
    $ws.Cells.Item($r, 2).Value = "Deliver Station"
}

# Update the selection to match the diff
$ws.Range("B2:B10").Select()
